$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 7-11, columns B:G with the new forecast error values

$ws.Range("B7").Value = 0.2458971023086367
$ws.Range("C7").Value = 1.638103339533811
$ws.Range("D7").Value = 4.276773671534449
$ws.Range("E7").Value = 2.068036187191716
$ws.Range("F7").Value = 2.082492313949006
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.3197323091405776
$ws.Range("C8").Value = 1.567098658972371
$ws.Range("D8").Value = 4.185152986416608
$ws.Range("E8").Value = 2.045764645900551
$ws.Range("F8").Value = 2.050124445948619
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 1.347797958024911
$ws.Range("C9").Value = 1.654972518212459
$ws.Range("D9").Value = 4.529824585228557
$ws.Range("E9").Value = 2.128338456455777
$ws.Range("F9").Value = 1.689990734814975
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = 0.976162992370311
$ws.Range("C10").Value = 1.401551600733488
$ws.Range("D10").Value = 3.042040935800302
$ws.Range("E10").Value = 1.744144757696534
$ws.Range("F10").Value = 1.504407627762786
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.8645624131408338
$ws.Range("C11").Value = 1.625693303335647
$ws.Range("D11").Value = 3.888939690769033
$ws.Range("E11").Value = 1.972039474952019
$ws.Range("F11").Value = 1.981625445358283
$ws.Range("G11").Value = 5
